$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Sync + Materi B.Indo": mark MATERI (D) / READY (E) checkmarks for the
# exam sessions that are now ready.
$ws.Range("D4").Value = "✔️"
$ws.Range("E4").Value = "✔️"
$ws.Range("D5").Value = "✔️"
$ws.Range("E5").Value = "✔️"
$ws.Range("D6").Value = "✔️"
$ws.Range("D14").Value = "✔️"
$ws.Range("D15").Value = "✔️"

# Restore the selection to the schedule table (A3:E17) as it was left in
# the authoring session.
[void]$ws.Range("A3:E17").Select()
